$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Valuation row (row 18): merge the "inside zipfolder" text into the
# URL cell, and mark the old standalone cell as not applicable.
$ws.Range("D18").Value = "~577projects\fall2013\projects\team09\team09a\Valuation(inside zipfolder(team09a.zip)"
$ws.Range("E18").Value = "NA"

# Remember the "file size" / model numbers for each project in the new column E.
$ws.Range("E6").Value = 541
$ws.Range("E7").Value = 577.6
$ws.Range("E8").Value = 912
$ws.Range("E9").Value = "NA"

$ws.Range("E10").Value = 1193
$ws.Range("E11").Value = 1328
$ws.Range("E12").Value = 512
$ws.Range("E13").Value = 432
$ws.Range("E15").Value = 384.56
$ws.Range("E16").Value = 1045.76
$ws.Range("E17").Value = 1810.32

$ws.Range("E19").Value = 1393.84
$ws.Range("E20").Value = 1092.8800000000001
$ws.Range("E21").Value = 8224.7199999999993
$ws.Range("E22").Value = 1741.92
$ws.Range("E23").Value = 950
$ws.Range("E24").Value = 1295.04
$ws.Range("E25").Value = 1506.32

# Leave the view scrolled back to the top and the whole row 5 selected, matching
# where the user was working.
$ws.Range("A5:XFD5").Select()
